$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells are treated as text (they store formatted
# numeric/percent strings, e.g. "306.96" or "6.52%") so Excel does not
# reinterpret them as numeric/percentage values.
$cells = @(
    "D2"
    "E2"
    "D3"
    "E3"
    "D4"
    "E4"
    "D5"
    "E5"
    "D6"
    "D7"
    "E7"
    "D8"
    "E8"
    "D9"
    "E9"
    "D10"
    "E10"
    "E11"
    "D12"
    "E12"
    "D13"
    "E13"
    "E14"
    "D15"
    "E15"
    "D16"
    "E16"
    "D17"
    "E17"
    "D18"
    "E18"
    "D19"
    "E19"
    "D20"
    "E20"
    "E21"
    "D22"
    "E22"
    "D23"
    "E23"
    "D24"
    "E24"
    "D25"
    "E25"
    "D26"
    "E26"
    "D27"
    "E27"
    "E28"
    "D40"
    "E40"
    "D41"
    "E41"
    "D42"
    "E42"
    "D43"
    "E43"
    "D44"
    "E44"
    "D45"
    "E45"
    "D46"
    "E46"
    "E47"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped for this run
$ws.Range("D2").Value = "306.96"
$ws.Range("E2").Value = "6.52%"
$ws.Range("D3").Value = "32.55"
$ws.Range("E3").Value = "10.87%"
$ws.Range("D4").Value = "5.293"
$ws.Range("E4").Value = "2.87%"
$ws.Range("D5").Value = "0.07434"
$ws.Range("E5").Value = "11.47%"
$ws.Range("D6").Value = "7.756"
$ws.Range("D7").Value = "3.711"
$ws.Range("E7").Value = "9.05%"
$ws.Range("D8").Value = "1.592"
$ws.Range("E8").Value = "16.89%"
$ws.Range("D9").Value = "0.9204"
$ws.Range("E9").Value = "0.54%"
$ws.Range("D10").Value = "0.01615"
$ws.Range("E10").Value = "2,389.77%"
$ws.Range("E11").Value = "5.93%"
$ws.Range("D12").Value = "0.07330"
$ws.Range("E12").Value = "11.67%"
$ws.Range("D13").Value = "0.07976"
$ws.Range("E13").Value = "3.52%"
$ws.Range("E14").Value = "5.52%"
$ws.Range("D15").Value = "0.09862"
$ws.Range("E15").Value = "9.66%"
$ws.Range("D16").Value = "0.001526"
$ws.Range("E16").Value = "-3.43%"
$ws.Range("D17").Value = "0.04558"
$ws.Range("E17").Value = "1.80%"
$ws.Range("D18").Value = "0.006169"
$ws.Range("E18").Value = "-1.18%"
$ws.Range("D19").Value = "3.480"
$ws.Range("E19").Value = "0.65%"
$ws.Range("D20").Value = "2.240"
$ws.Range("E20").Value = "0.81%"
$ws.Range("E21").Value = "1.87%"
$ws.Range("D22").Value = "0.1320"
$ws.Range("E22").Value = "0.84%"
$ws.Range("D23").Value = "4.245"
$ws.Range("E23").Value = "4.61%"
$ws.Range("D24").Value = "0.1619"
$ws.Range("E24").Value = "4.40%"
$ws.Range("D25").Value = "0.001195"
$ws.Range("E25").Value = "0.45%"
$ws.Range("D26").Value = "0.004538"
$ws.Range("E26").Value = "9.77%"
$ws.Range("D27").Value = "0.0001168"
$ws.Range("E27").Value = "-6.43%"
$ws.Range("E28").Value = "2.94%"
$ws.Range("D40").Value = "0.04499"
$ws.Range("E40").Value = "6.89%"
$ws.Range("D41").Value = "0.007297"
$ws.Range("E41").Value = "8.51%"
$ws.Range("D42").Value = "0.1365"
$ws.Range("E42").Value = "10.01%"
$ws.Range("D43").Value = "0.002256"
$ws.Range("E43").Value = "14.10%"
$ws.Range("D44").Value = "0.01368"
$ws.Range("E44").Value = "6.88%"
$ws.Range("D45").Value = "0.00005969"
$ws.Range("E45").Value = "6.74%"
$ws.Range("D46").Value = "1.892"
$ws.Range("E46").Value = "-3.83%"
$ws.Range("E47").Value = "-0.54%"

# Restore default (General) styling so no stray text-format style lingers
foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Updated $($cells.Count) cells"
